$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) rows 2-41: update date serial from 45701 (2025-02-13)
# to 45702 (2025-02-14) while keeping the existing date formatting/style.
for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45701) {
        $cell.Value = 45702
    }
}
